$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First weighting table (rows 27-32, adjustments block) ---
$ws.Range("E27").Value = "weight"
$ws.Range("F27").Value = "weighted x"
$ws.Range("G27").Value = "weighted y"

$ws.Range("E28").Value = 1
$ws.Range("F28").Formula = "=C28*E28"
$ws.Range("G28").Formula = "=E28*D28"

$ws.Range("E29").Value = 1
$ws.Range("E30").Value = 1
$ws.Range("E31").Value = 1

$ws.Range("F29:F31").Formula = "=C29*E29"
$ws.Range("G29:G31").Formula = "=E29*D29"

$ws.Range("F32").Formula = "=SUM(F28:F31)"
$ws.Range("G32").Formula = "=SUM(G28:G31)"

# --- Second weighting table (rows 36-42, velocity block) ---
$ws.Range("E36").Value = "weight"
$ws.Range("F36").Value = "weighted"

$ws.Range("E37").Value = 1
$ws.Range("F37").Formula = "=C37*E37"
$ws.Range("G37").Formula = "=D37*E37"

$ws.Range("E38").Value = 1
$ws.Range("E39").Value = 1
$ws.Range("E40").Value = 0.1
$ws.Range("E41").Value = 1

$ws.Range("F38:F41").Formula = "=C38*E38"
$ws.Range("G38:G41").Formula = "=D38*E38"

$ws.Range("F42").Formula = "=SUM(F37:F41)"
$ws.Range("G42").Formula = "=SUM(G37:G41)"

# --- Final selection / view state ---
$ws.Range("G42").Select() | Out-Null
